$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: CONVERSION
# ---------------------------------------------------------------
$wsConv = $wb.Worksheets.Item("CONVERSION")

# Remove the now-unused deliveryTemp_degc / capacityHeat_kW columns (F:G)
$wsConv.Range("F1:G3").EntireColumn.Delete()

# Row 2 becomes the new DIESEL_VEHICLE entry
$wsConv.Range("B2").Value = "DIESEL_VEHICLE"
$wsConv.Range("C2").Value = "VehicleConversionAsset"
$wsConv.Range("D2").Value = "Diesel_Truck"
$wsConv.Range("E2").Value = 0.95

# Row 3 (GAS_BURNER) stays as-is (category/type/dataclassname/name/eta_r unchanged)

# New row 4: duplicate of the original GAS_BURNER entry
$wsConv.Range("A4").Value = "CONVERSION"
$wsConv.Range("B4").Value = "GAS_BURNER"
$wsConv.Range("C4").Value = "ChemicalHeatConversionAsset"
$wsConv.Range("D4").Value = "Building_gas_burner_60kW"
$wsConv.Range("E4").Value = 0.95

# ---------------------------------------------------------------
# Sheet: CONSUMPTION
# ---------------------------------------------------------------
$wsCons = $wb.Worksheets.Item("CONSUMPTION")

$wsCons.Range("E1").Value = "yearlyDemandHeat_kWh"

$wsCons.Range("B2").Value = "HEAT_DEMAND"
$wsCons.Range("C2").Value = "HeatConsumptionAsset"
$wsCons.Range("D2").Value = "INDUSTRY_OTHER_HEAT_DEMAND"
$wsCons.Range("E2").Value = 600000

# ---------------------------------------------------------------
# Sheet: PRODUCTION
# ---------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("PRODUCTION")

$wsProd.Range("D3").Value = "Building_solarpanels_10kWp"
$wsProd.Range("E3").Value = 10

$wsProd.Range("D4").Value = "Building_solarpanels_0kWp"
$wsProd.Range("E4").Value = 0

# New row 6: duplicate of the Solarpanels_1MW entry
$wsProd.Range("A6").Value = "PRODUCTION"
$wsProd.Range("B6").Value = "PHOTOVOLTAIC"
$wsProd.Range("C6").Value = "ElectricProductionAsset"
$wsProd.Range("D6").Value = "Solarpanels_1MW"
$wsProd.Range("E6").Value = 1000

# ---------------------------------------------------------------
# Sheet: STORAGE
# ---------------------------------------------------------------
$wsStor = $wb.Worksheets.Item("STORAGE")

# Reorder header columns: D/E/F now hold stateOfCharge_r / name / capacityElectricity_kW
$wsStor.Range("D1").Value = "stateOfCharge_r"
$wsStor.Range("E1").Value = "name"
$wsStor.Range("F1").Value = "capacityElectricity_kW"
# G1 (storageCapacity_kWh) is unchanged

# The previous single STORAGE_ELECTRIC / Grid_battery_10MWh row is now pushed to row 9,
# with columns re-mapped to the new layout.
$wsStor.Range("A9").Value = "STORAGE"
$wsStor.Range("B9").Value = "STORAGE_ELECTRIC"
$wsStor.Range("C9").Value = "ElectricStorageAsset"
$wsStor.Range("D9").Value = 1
$wsStor.Range("E9").Value = "Grid_battery_10MWh"
$wsStor.Range("F9").Value = 2000
$wsStor.Range("G9").Value = 10000

# New row 8: STORAGE_ELECTRIC / Grid_battery_7MWh
$wsStor.Range("A8").Value = "STORAGE"
$wsStor.Range("B8").Value = "STORAGE_ELECTRIC"
$wsStor.Range("C8").Value = "ElectricStorageAsset"
$wsStor.Range("D8").Value = 1
$wsStor.Range("E8").Value = "Grid_battery_7MWh"
$wsStor.Range("F8").Value = 1000
$wsStor.Range("G8").Value = 7000

# New rows 2-7: six ELECTRIC_HEAVY_GOODS_VEHICLE entries
for ($r = 2; $r -le 7; $r++) {
    $wsStor.Range("A$r").Value = "STORAGE"
    $wsStor.Range("B$r").Value = "ELECTRIC_HEAVY_GOODS_VEHICLE"
    $wsStor.Range("C$r").Value = "VehicleElectricStorageAsset"
    $wsStor.Range("D$r").Value = 1
    $wsStor.Range("E$r").Value = "EHGV"
    $wsStor.Range("F$r").Value = 100
    $wsStor.Range("G$r").Value = 500
}
